$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 5) to the table:
#   Tên nhóm | Thời gian hoàn thành | Chức năng | Thành viên hoặc nhóm phát triển
$ws.Range("A5").Value = "The Owls"
$ws.Range("B5").Value = "21/12/2017"
# Set D5 before C5 so the shared-string table records "Trần Nguyên (1412360)"
# ahead of "Quán lý thuê phòng (xem, sửa, thêm)".
$ws.Range("D5").Value = "Trần Nguyên (1412360)"
$ws.Range("C5").Value = "Quán lý thuê phòng (xem, sửa, thêm)"

# Move the active selection the way Excel would after typing the new row.
[void]$ws.Range("B6").Select()
